$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Bro" column header
$ws.Range("F1").Value = "Bro"

# Clear a handful of individual "missing values" cells
$ws.Range("E4").ClearContents()
$ws.Range("B7").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("E10").ClearContents()

# Populate the new "Bro" column with Yes/No values
$ws.Range("F2").Value = "Yes"
$ws.Range("F3").Value = "Yes"
$ws.Range("F4").Value = "Yes"
$ws.Range("F5").Value = "No"
$ws.Range("F6").Value = "No"
$ws.Range("F7").Value = "Yes"
$ws.Range("F8").Value = "No"
$ws.Range("F10").Value = "Yes"
$ws.Range("F11").Value = "Yes"

# Update the selection to match the saved workbook state
$ws.Range("E10").Select()
